# Fruta / hortaliza, semanal
# Re-order the weekly price records (rows 2-17) according to the new
# upstream data pull. Only columns D (Fecha), M (Volumen), N (Precio
# minimo), O (Precio maximo), P (Precio promedio ponderado), R (Origen)
# and S (Precio $/Kg) move between rows; everything else is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values for the columns that get reshuffled,
# keyed by their original row number.
$before = @{}
for ($r = 2; $r -le 17; $r++) {
    $before[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        R = $ws.Cells.Item($r, 18).Value2
        S = $ws.Cells.Item($r, 19).Value2
    }
}

# Maps the destination row -> source row (which "before" row's data now
# lands on this row).
$rowMap = @{
    2  = 5
    3  = 12
    4  = 6
    5  = 10
    6  = 3
    7  = 9
    8  = 8
    9  = 16
    10 = 15
    11 = 14
    12 = 2
    13 = 17
    14 = 7
    15 = 13
    16 = 11
    17 = 4
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $vals = $before[$srcRow]

    $ws.Cells.Item($destRow, 4).Value2  = $vals.D
    $ws.Cells.Item($destRow, 13).Value2 = $vals.M
    $ws.Cells.Item($destRow, 14).Value2 = $vals.N
    $ws.Cells.Item($destRow, 15).Value2 = $vals.O
    $ws.Cells.Item($destRow, 16).Value2 = $vals.P
    $ws.Cells.Item($destRow, 18).Value2 = $vals.R
    $ws.Cells.Item($destRow, 19).Value2 = $vals.S
}
